$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert "Republica de Africa Central" right after "Nicaragua" (and before "Georgia") ---
# This pushes Georgia / Jordania / Republica del Chad / Crucero down by one row each.
$ws.Range("A127").Value = "Republica de Africa Central"
$ws.Range("A128").Value = "Georgia"
$ws.Range("A129").Value = "Jordania"
$ws.Range("A130").Value = "Republica del Chad"
$ws.Range("A131").Value = "Crucero"

# --- Refresh the "last updated" banner ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Mayo de 2020 a las 00:40"

# --- Refresh COVID-19 numbers for Estados Unidos (row 4) ---
$ws.Range("B4").Value = 1766795
$ws.Range("C4").Value = 20992
$ws.Range("D4").Value = 498013
$ws.Range("E4").Value = 1165486
$ws.Range("G4").Value = 1189
$ws.Range("H4").Value = 103296

# --- Peru (row 15) ---
$ws.Range("B15").Value = 141779
$ws.Range("C15").Value = 5874
$ws.Range("D15").Value = 59442
$ws.Range("E15").Value = 78238
$ws.Range("G15").Value = 116
$ws.Range("H15").Value = 4099

# --- Canada (row 16) ---
$ws.Range("B16").Value = 88504
$ws.Range("C16").Value = 985
$ws.Range("D16").Value = 46831
$ws.Range("E16").Value = 34798
$ws.Range("G16").Value = 110
$ws.Range("H16").Value = 6875

# --- Sudafrica (row 33) ---
$ws.Range("B33").Value = 27403
$ws.Range("C33").Value = 1466
$ws.Range("D33").Value = 14370
$ws.Range("E33").Value = 12456
$ws.Range("G33").Value = 25
$ws.Range("H33").Value = 577

# --- Japon (row 43) ---
$ws.Range("B43").Value = 16683
$ws.Range("C43").Value = 32
$ws.Range("D43").Value = 14147
$ws.Range("E43").Value = 1669
$ws.Range("G43").Value = 9
$ws.Range("H43").Value = 867

# --- Noruega (row 59) ---
$ws.Range("B59").Value = 8411
$ws.Range("C59").Value = 10
$ws.Range("E59").Value = 448

# --- Gabon (row 87) ---
$ws.Range("B87").Value = 2431
$ws.Range("C87").Value = 112
$ws.Range("D87").Value = 668
$ws.Range("E87").Value = 1749

# --- Republica de Africa Central (now row 127) ---
$ws.Range("B127").Value = 755
$ws.Range("C127").Value = 53
$ws.Range("D127").Value = 23
$ws.Range("E127").Value = 731
$ws.Range("H127").Value = 1

# --- Georgia (now row 128) ---
$ws.Range("B128").Value = 738
$ws.Range("C128").Value = 3
$ws.Range("D128").Value = 573
$ws.Range("E128").Value = 153
$ws.Range("H128").Value = 12

# --- Jordania (now row 129) ---
$ws.Range("B129").Value = 728
$ws.Range("C129").Value = 8
$ws.Range("D129").Value = 497
$ws.Range("E129").Value = 222
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 9

# --- Republica del Chad (now row 130) ---
$ws.Range("B130").Value = 726
$ws.Range("C130").Value = 11
$ws.Range("D130").Value = 413
$ws.Range("E130").Value = 248
$ws.Range("G130").Value = 1
$ws.Range("H130").Value = 65

# --- Crucero (now row 131) ---
$ws.Range("B131").Value = 712
$ws.Range("D131").Value = 651
$ws.Range("E131").Value = 48
$ws.Range("H131").Value = 13

# --- Guayana Francesa (row 143) ---
$ws.Range("B143").Value = 409
$ws.Range("C143").Value = 3
$ws.Range("E143").Value = 258
